$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value2 = "CATEGORÍA"
$ws.Range("B1").Value2 = "SUBCATEGORÍA"
$ws.Range("C1").Value2 = "DESCRIPCIÓN"
$ws.Range("D1").Value2 = "FECHA"

# Row 2
$ws.Range("A2").Value2 = "Procesos internauticos"
$ws.Range("B2").Value2 = "Interpolación de cadenas concatenadas"
$ws.Range("C2").Value2 = "Esta tarea es muy macabra"
$ws.Range("D2").Value2 = "2021-05-26 11:23:46"

# Row 3
$ws.Range("A3").Value2 = "Procesos internauticos"
$ws.Range("B3").Value2 = "Maravillas modernas de los servidores publicos"
$ws.Range("C3").Value2 = "Aquí inserto el texto de la tarea`nbla bla bla"
$ws.Range("D3").Value2 = "2021-05-28 12:21:48"

# Row 4
$ws.Range("A4").Value2 = "Pansensuales"
$ws.Range("B4").Value2 = "Amo el pan"
$ws.Range("C4").Value2 = "Esta descripción está super detallada ya que presenta sinopsis sinapticas en las mamalasticas"
$ws.Range("D4").Value2 = "2021-05-26 12:31:52"

# Row 5
$ws.Range("A5").Value2 = "Pansensuales"
$ws.Range("B5").Value2 = "Amo el pan"
$ws.Range("C5").Value2 = "Otra tarea dista mondá"
$ws.Range("D5").Value2 = "2021-05-28 12:40:21"

# Row 6
$ws.Range("A6").Value2 = "Procesos internauticos"
$ws.Range("B6").Value2 = "Maravillas modernas de los servidores publicos"
$ws.Range("C6").Value2 = "Esta será la sexta tarea"
$ws.Range("D6").Value2 = "2021-05-26 17:34:19"

# Row 7
$ws.Range("A7").Value2 = "Pansensuales"
$ws.Range("B7").Value2 = "Amo la mogolla"
$ws.Range("C7").Value2 = "Tarea super genérica"
$ws.Range("D7").Value2 = "2021-05-28 09:27:10"

# Row 8
$ws.Range("A8").Value2 = "Pansensuales"
$ws.Range("B8").Value2 = "Amo el pan"
$ws.Range("C8").Value2 = "BASADO recontramamabuebo"
$ws.Range("D8").Value2 = "2021-05-31 23:08:01"

# Row 9
$ws.Range("A9").Value2 = "Phva pic"
$ws.Range("B9").Value2 = "Contratación"
$ws.Range("C9").Value2 = "se realizó el documento. `nMama guevo"
$ws.Range("D9").Value2 = "2021-06-02 18:15:02"

# Row 10
$ws.Range("A10").Value2 = "Phva pic"
$ws.Range("B10").Value2 = "Seguimiento"
$ws.Range("C10").Value2 = "a ver"
$ws.Range("D10").Value2 = "2021-06-23 20:30:52"

# Row 11
$ws.Range("A11").Value2 = "catejemplo 1"
$ws.Range("B11").Value2 = "a ver que pasa"
$ws.Range("C11").Value2 = "a ver que pasa"
$ws.Range("D11").Value2 = "2021-06-23 20:33:32"

# Row 12
$ws.Range("A12").Value2 = "1. phva pic"
$ws.Range("B12").Value2 = "1. ejecución interno o alianzas"
$ws.Range("C12").Value2 = "Prueba de mamarre"
$ws.Range("D12").Value2 = "2021-08-08 18:37:57"

# Remove now-obsolete rows 13 and 14 (report was trimmed to 12 rows)
$ws.Rows("13:14").Delete() | Out-Null

# Restore standard row height for rows whose multi-line text triggered autofit
$ws.Rows(3).AutoFit() | Out-Null
$ws.Rows(9).AutoFit() | Out-Null
